$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 32.120528
$ws.Range("H2").Value = 96.36158399999999
$ws.Range("I2").Value = 0.5492202673743384
$ws.Range("J2").Value = 0.5492202673743384
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.52818033333334
$ws.Range("N2").Value = 52.584541
$ws.Range("O2").Value = 0.5698985076516194
$ws.Range("P2").Value = 0.5698985076516194
$ws.Range("Q2").Value = 563.0144071858828
$ws.Range("R2").Value = 5067.129664672943
$ws.Range("S2").Value = 0.3129998107486588
$ws.Range("T2").Value = 0.3129998107486588
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 32.120528
$ws.Range("H3").Value = 96.36158399999999
$ws.Range("I3").Value = 0.5492202673743384
$ws.Range("J3").Value = 0.5492202673743384
$ws.Range("O3").Value = 0.3813545701360604
$ws.Range("P3").Value = 0.3813545701360604
$ws.Range("Q3").Value = 376.7479899491734
$ws.Range("R3").Value = 3390.73190954256
$ws.Range("S3").Value = 0.209447658974553
$ws.Range("T3").Value = 0.209447658974553
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 32.120528
$ws.Range("H4").Value = 96.36158399999999
$ws.Range("I4").Value = 0.5492202673743384
$ws.Range("J4").Value = 0.5492202673743384
$ws.Range("M4").Value = 1.305498666666667
$ws.Range("N4").Value = 3.916496
$ws.Range("O4").Value = 0.04244603419897754
$ws.Range("P4").Value = 0.04244603419897754
$ws.Range("Q4").Value = 41.93330647662933
$ws.Range("R4").Value = 377.399758289664
$ws.Range("S4").Value = 0.02331222225174276
$ws.Range("T4").Value = 0.02331222225174276
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 32.120528
$ws.Range("H5").Value = 96.36158399999999
$ws.Range("I5").Value = 0.5492202673743384
$ws.Range("J5").Value = 0.5492202673743384
$ws.Range("M5").Value = 0.1937943333333333
$ws.Range("N5").Value = 0.581383
$ws.Range("O5").Value = 0.00630088801334258
$ws.Range("P5").Value = 0.00630088801334258
$ws.Range("Q5").Value = 6.224776310074666
$ws.Range("R5").Value = 56.02298679067199
$ws.Range("S5").Value = 0.003460575399383776
$ws.Range("T5").Value = 0.003460575399383776
$ws.Range("G6").Value = 22.49179166666667
$ws.Range("H6").Value = 67.475375
$ws.Range("I6").Value = 0.3845810950833243
$ws.Range("J6").Value = 0.3845810950833244
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.52818033333334
$ws.Range("N6").Value = 52.584541
$ws.Range("O6").Value = 0.5698985076516194
$ws.Range("P6").Value = 0.5698985076516194
$ws.Range("Q6").Value = 394.2401803530973
$ws.Range("R6").Value = 3548.161623177875
$ws.Range("S6").Value = 0.2191721921590121
$ws.Range("T6").Value = 0.2191721921590121
$ws.Range("G7").Value = 22.49179166666667
$ws.Range("H7").Value = 67.475375
$ws.Range("I7").Value = 0.3845810950833243
$ws.Range("J7").Value = 0.3845810950833244
$ws.Range("O7").Value = 0.3813545701360604
$ws.Range("P7").Value = 0.3813545701360604
$ws.Range("Q7").Value = 263.8106478440278
$ws.Range("R7").Value = 2374.29583059625
$ws.Range("S7").Value = 0.1466617581979565
$ws.Range("T7").Value = 0.1466617581979566
$ws.Range("G8").Value = 22.49179166666667
$ws.Range("H8").Value = 67.475375
$ws.Range("I8").Value = 0.3845810950833243
$ws.Range("J8").Value = 0.3845810950833244
$ws.Range("M8").Value = 1.305498666666667
$ws.Range("N8").Value = 3.916496
$ws.Range("O8").Value = 0.04244603419897754
$ws.Range("P8").Value = 0.04244603419897754
$ws.Range("Q8").Value = 29.36300403177778
$ws.Range("R8").Value = 264.267036286
$ws.Range("S8").Value = 0.01632394231418702
$ws.Range("T8").Value = 0.01632394231418702
$ws.Range("G9").Value = 22.49179166666667
$ws.Range("H9").Value = 67.475375
$ws.Range("I9").Value = 0.3845810950833243
$ws.Range("J9").Value = 0.3845810950833244
$ws.Range("M9").Value = 0.1937943333333333
$ws.Range("N9").Value = 0.581383
$ws.Range("O9").Value = 0.00630088801334258
$ws.Range("P9").Value = 0.00630088801334258
$ws.Range("Q9").Value = 4.358781771513889
$ws.Range("R9").Value = 39.229035943625
$ws.Range("S9").Value = 0.002423202412168681
$ws.Range("T9").Value = 0.002423202412168682
$ws.Range("G10").Value = 3.864125333333333
$ws.Range("H10").Value = 11.592376
$ws.Range("I10").Value = 0.06607163956773338
$ws.Range("J10").Value = 0.06607163956773338
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 17.52818033333334
$ws.Range("N10").Value = 52.584541
$ws.Range("O10").Value = 0.5698985076516194
$ws.Range("P10").Value = 0.5698985076516194
$ws.Range("Q10").Value = 67.73108567326845
$ws.Range("R10").Value = 609.579771059416
$ws.Range("S10").Value = 0.03765412878774694
$ws.Range("T10").Value = 0.03765412878774694
$ws.Range("G11").Value = 3.864125333333333
$ws.Range("H11").Value = 11.592376
$ws.Range("I11").Value = 0.06607163956773338
$ws.Range("J11").Value = 0.06607163956773338
$ws.Range("O11").Value = 0.3813545701360604
$ws.Range("P11").Value = 0.3813545701360604
$ws.Range("Q11").Value = 45.32308597931556
$ws.Range("R11").Value = 407.90777381384
$ws.Range("S11").Value = 0.02519672170553769
$ws.Range("T11").Value = 0.02519672170553769
$ws.Range("G12").Value = 3.864125333333333
$ws.Range("H12").Value = 11.592376
$ws.Range("I12").Value = 0.06607163956773338
$ws.Range("J12").Value = 0.06607163956773338
$ws.Range("M12").Value = 1.305498666666667
$ws.Range("N12").Value = 3.916496
$ws.Range("O12").Value = 0.04244603419897754
$ws.Range("P12").Value = 0.04244603419897754
$ws.Range("Q12").Value = 5.044610470499555
$ws.Range("R12").Value = 45.401494234496
$ws.Range("S12").Value = 0.002804479072674529
$ws.Range("T12").Value = 0.002804479072674529
$ws.Range("G13").Value = 3.864125333333333
$ws.Range("H13").Value = 11.592376
$ws.Range("I13").Value = 0.06607163956773338
$ws.Range("J13").Value = 0.06607163956773338
$ws.Range("M13").Value = 0.1937943333333333
$ws.Range("N13").Value = 0.581383
$ws.Range("O13").Value = 0.00630088801334258
$ws.Range("P13").Value = 0.00630088801334258
$ws.Range("Q13").Value = 0.7488455928897777
$ws.Range("R13").Value = 6.739610336008
$ws.Range("S13").Value = 0.0004163100017742226
$ws.Range("T13").Value = 0.0004163100017742226
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.007427333333333334
$ws.Range("H14").Value = 0.022282
$ws.Range("I14").Value = 0.0001269979746040186
$ws.Range("J14").Value = 0.0001269979746040186
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 17.52818033333334
$ws.Range("N14").Value = 52.584541
$ws.Range("O14").Value = 0.5698985076516194
$ws.Range("P14").Value = 0.5698985076516194
$ws.Range("Q14").Value = 0.1301876380624445
$ws.Range("R14").Value = 1.171688742562
$ws.Range("S14").Value = 0.00007237595620160848
$ws.Range("T14").Value = 0.00007237595620160849
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.007427333333333334
$ws.Range("H15").Value = 0.022282
$ws.Range("I15").Value = 0.0001269979746040186
$ws.Range("J15").Value = 0.0001269979746040186
$ws.Range("O15").Value = 0.3813545701360604
$ws.Range("P15").Value = 0.3813545701360604
$ws.Range("Q15").Value = 0.08711665337555556
$ws.Range("R15").Value = 0.78404988038
$ws.Range("S15").Value = 0.00004843125801326584
$ws.Range("T15").Value = 0.00004843125801326585
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.007427333333333334
$ws.Range("H16").Value = 0.022282
$ws.Range("I16").Value = 0.0001269979746040186
$ws.Range("J16").Value = 0.0001269979746040186
$ws.Range("M16").Value = 1.305498666666667
$ws.Range("N16").Value = 3.916496
$ws.Range("O16").Value = 0.04244603419897754
$ws.Range("P16").Value = 0.04244603419897754
$ws.Range("Q16").Value = 0.009696373763555555
$ws.Range("R16").Value = 0.087267363872
$ws.Range("S16").Value = 0.000005390560373243056
$ws.Range("T16").Value = 0.000005390560373243057
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.007427333333333334
$ws.Range("H17").Value = 0.022282
$ws.Range("I17").Value = 0.0001269979746040186
$ws.Range("J17").Value = 0.0001269979746040186
$ws.Range("M17").Value = 0.1937943333333333
$ws.Range("N17").Value = 0.581383
$ws.Range("O17").Value = 0.00630088801334258
$ws.Range("P17").Value = 0.00630088801334258
$ws.Range("Q17").Value = 0.001439375111777778
$ws.Range("R17").Value = 0.012954376006
$ws.Range("S17").Value = 0.0000008002000159012463
$ws.Range("T17").Value = 0.0000008002000159012464
